$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 232, shifting existing rows 232:360 down to 233:361
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new data record
$ws.Cells.Item(232, 1).Value = 10
$ws.Cells.Item(232, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(232, 3).Value = "La Araucanía"
$ws.Cells.Item(232, 4).Value = 44806
$ws.Cells.Item(232, 5).Value = 9
$ws.Cells.Item(232, 6).Value = 100112017
$ws.Cells.Item(232, 7).Value = "Apio"
$ws.Cells.Item(232, 8).Value = "Americana (o)"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 115
$ws.Cells.Item(232, 11).Value = 9000
$ws.Cells.Item(232, 12).Value = 10000
$ws.Cells.Item(232, 13).Value = 9565
$ws.Cells.Item(232, 14).Value = "$/docena de matas"
$ws.Cells.Item(232, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(232, 16).Value = 1594
$ws.Cells.Item(232, 17).Value = 6
$ws.Cells.Item(232, 18).Value = "Hortaliza"
